$wb = $excel.ActiveWorkbook

# --- Sheet "Means": update row 9 (Total Cancer Risk per million) and row 10 (Total Respiratory hazard quotient)
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 27
$wsMeans.Range("D9").Value = 30
$wsMeans.Range("E9").Value = 36
$wsMeans.Range("F9").Value = 33
$wsMeans.Range("G9").Value = 30

$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.33
$wsMeans.Range("D10").Value = 0.3
$wsMeans.Range("E10").Value = 0.3
$wsMeans.Range("F10").Value = 0.3
$wsMeans.Range("G10").Value = 0.3

# --- Sheet "Standard Deviations": update row 9 and row 10
$wsStd = $wb.Worksheets.Item("Standard Deviations")

$wsStd.Range("B9").Value = 8.3
$wsStd.Range("C9").Value = 6.4
$wsStd.Range("D9").Value = 0
$wsStd.Range("E9").Value = 15
$wsStd.Range("F9").Value = 9.1
$wsStd.Range("G9").Value = 5.2

$wsStd.Range("B10").Value = 0.11
$wsStd.Range("C10").Value = 0.07
$wsStd.Range("D10").Value = 0
$wsStd.Range("E10").Value = 0
$wsStd.Range("F10").Value = 0
$wsStd.Range("G10").Value = 0.011

$wb.Save()
